# semana 40 de 2025
# Add week 39 (AP) and week 40 (AQ) columns to the weekly IRA extract.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new week-number labels, matching the text style
# already used by the other week headers (bold + centered, stored as text
# rather than numbers - temporarily mark the cells as Text so "39"/"40"
# aren't auto-coerced into numeric values, then restore General format).
$ap1 = $ws.Cells.Item(1, 42)
$aq1 = $ws.Cells.Item(1, 43)

$ap1.NumberFormat = "@"
$ap1.Value = "39"
$ap1.NumberFormat = "General"

$aq1.NumberFormat = "@"
$aq1.Value = "40"
$aq1.NumberFormat = "General"

$headerRange = $ws.Range("AP1:AQ1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108

# --- Data rows: counts for week 39 (AP) and week 40 (AQ) per establishment.
# Rows not listed here (4, 18-21, 27, 32-34, 39, 40) had no reported cases
# for either week and are intentionally left blank, matching the source.
$weekData = @{
    2  = @(54, 66)
    3  = @(77, 68)
    5  = @(0, 2)
    6  = @(122, 109)
    7  = @(30, 29)
    8  = @(39, 32)
    9  = @(3, 3)
    10 = @(2, 2)
    11 = @(1, $null)
    12 = @(2, 3)
    13 = @(1, $null)
    14 = @(4, 3)
    15 = @($null, 2)
    16 = @(1, 1)
    17 = @(2, 1)
    22 = @(1, 1)
    23 = @(7, 5)
    24 = @($null, 3)
    25 = @(49, 39)
    26 = @(3, $null)
    28 = @(229, 218)
    29 = @(0, 0)
    30 = @(84, 79)
    31 = @(2, 3)
    35 = @(55, 48)
    36 = @(1, 0)
    37 = @(17, 10)
    38 = @(74, 82)
    41 = @(13, 9)
    42 = @(8, 6)
    43 = @(26, 25)
    44 = @(177, $null)
    45 = @(95, 99)
    46 = @(143, 182)
    47 = @(1, 2)
    48 = @(107, 148)
    49 = @(4, 6)
    50 = @(0, 0)
    51 = @(3, 9)
    52 = @(3, 3)
    53 = @(7, 15)
    54 = @(0, 0)
    55 = @(0, 0)
    56 = @(2, 6)
    57 = @(121, 31)
    58 = @(12, 14)
}

foreach ($row in $weekData.Keys) {
    $pair = $weekData[$row]
    $ap = $pair[0]
    $aq = $pair[1]
    if ($null -ne $ap) {
        $ws.Cells.Item($row, 42).Value = $ap
    }
    if ($null -ne $aq) {
        $ws.Cells.Item($row, 43).Value = $aq
    }
}
